# Weekly fruit/hortaliza data update:
# Insert two new rows (new week of data) right before the current row 67,
# pushing all existing data rows (old 67..98) down by two rows to 69..100.
# Then populate the two newly inserted rows (67 and 68) with the new week's
# price data (Primera / Segunda quality for "Locoto").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(67).Insert()
$ws.Rows.Item(67).Insert()

# New row 67 - "Primera" quality
$ws.Cells.Item(67, 1).Value = 1
$ws.Cells.Item(67, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(67, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(67, 4).Value = 44603
$ws.Cells.Item(67, 5).Value = 15
$ws.Cells.Item(67, 6).Value = 100112042
$ws.Cells.Item(67, 7).Value = "Locoto"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 130
$ws.Cells.Item(67, 11).Value = 33000
$ws.Cells.Item(67, 12).Value = 35000
$ws.Cells.Item(67, 13).Value = 34000
$ws.Cells.Item(67, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(67, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(67, 16).Value = 1700
$ws.Cells.Item(67, 17).Value = 20
$ws.Cells.Item(67, 18).Value = "Hortaliza"

# New row 68 - "Segunda" quality
$ws.Cells.Item(68, 1).Value = 1
$ws.Cells.Item(68, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value = 44603
$ws.Cells.Item(68, 5).Value = 15
$ws.Cells.Item(68, 6).Value = 100112042
$ws.Cells.Item(68, 7).Value = "Locoto"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Segunda"
$ws.Cells.Item(68, 10).Value = 150
$ws.Cells.Item(68, 11).Value = 28000
$ws.Cells.Item(68, 12).Value = 29000
$ws.Cells.Item(68, 13).Value = 28500
$ws.Cells.Item(68, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 1425
$ws.Cells.Item(68, 17).Value = 20
$ws.Cells.Item(68, 18).Value = "Hortaliza"
